$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 values (C1:F1), clear G1 entirely
$ws.Range("C1").Value = -0.278
$ws.Range("D1").Value = -122.4354
$ws.Range("E1").Value = 42.8765
$ws.Range("F1").Value = 179.7921
$ws.Range("G1").ClearContents()

# Row 1 height matches the rest of the sheet now
$ws.Rows.Item(1).RowHeight = 13.8

# Row 2 is fully cleared (was A2:G2)
$ws.Range("A2:G2").ClearContents()

# Update the active selection to H5, as in the target workbook
$ws.Range("H5").Select()
